$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "ángel perdomo",
    "brent honeywell",
    "carl edwards",
    "daniel lynch",
    "duane underwood",
    "j.b. bukauskas",
    "j.t. chargois",
    "jaime barría",
    "jose ferrer",
    "julio teherán",
    "mark leiter",
    "matt boyd",
    "matthew bowman",
    "mike king",
    "néstor cortés",
    "nick martínez"
)

$eras = @(3.72, 4.82, 3.69, 4.64, 5.18, 1.29, 3.61, 5.68, 5.03, 4.4, 3.5, 5.45, 9, 2.75, 4.97, 3.43)

$startRow = 865
$startIndex = 863

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i
    # Copy formatting (bold/border/centered style) from the last existing
    # data row in column A before writing the new value into it.
    $ws.Range("A864").Copy($ws.Range("A$row"))
    $ws.Cells.Item($row, 1).Value = $startIndex + $i
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $eras[$i]
}
